$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 30 (new rows 31 and 32), shifting existing
# rows 31.. down by two.
$ws.Rows("31:32").Insert()

# New row 31 data
$ws.Range("A31").Value = 10
$ws.Range("B31").Value = "Vega Modelo de Temuco"
$ws.Range("C31").Value = "La Araucanía"
$ws.Range("D31").Value = 44623
$ws.Range("E31").Value = 9
$ws.Range("F31").Value = 100112031
$ws.Range("G31").Value = "Poroto verde"
$ws.Range("H31").Value = "Brío"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 80
$ws.Range("K31").Value = 1200
$ws.Range("L31").Value = 1200
$ws.Range("M31").Value = 1200
$ws.Range("N31").Value = "`$/kilo"
$ws.Range("O31").Value = "Región de La Araucanía"
$ws.Range("P31").Value = 1200
$ws.Range("Q31").Value = 1
$ws.Range("R31").Value = "Hortaliza"

# New row 32 data
$ws.Range("A32").Value = 10
$ws.Range("B32").Value = "Vega Modelo de Temuco"
$ws.Range("C32").Value = "La Araucanía"
$ws.Range("D32").Value = 44623
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = 100112031
$ws.Range("G32").Value = "Poroto verde"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 80
$ws.Range("K32").Value = 1200
$ws.Range("L32").Value = 1200
$ws.Range("M32").Value = 1200
$ws.Range("N32").Value = "`$/kilo"
$ws.Range("O32").Value = "Región de La Araucanía"
$ws.Range("P32").Value = 1200
$ws.Range("Q32").Value = 1
$ws.Range("R32").Value = "Hortaliza"
